# Marksheet update: recompute the score summary (rows 10-12) using numeric
# values (instead of leftover inline-string artifacts), refresh the
# "Correct answer" column so it highlights right/wrong answers with the
# correct/incorrect styles, and drop the now-unused extra
# "Student Ans/Correct Ans" column blocks (columns D:E past row 18 and the
# whole G:H block) that used to hold duplicate attempt data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the "mtitleStyle" (bold, centered) look used by the other
#        row headers (No. / Marking / Total) to A10:A12. We copy formats
#        only from A9, which already carries that exact style, so we reuse
#        the existing style index instead of minting a new one. ---
foreach ($addr in @("A10", "A11", "A12")) {
    $ws.Range("A9").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- 2. Apply the "correctStyle" (green, centered) look to the cells that
#        now flag a correct answer: D16 and the "A" column entries for
#        questions 21, 23, 28 and 32. Copy formats only from B10, which
#        already carries that exact style. ---
foreach ($addr in @("D16", "A21", "A23", "A28", "A32")) {
    $ws.Range("B10").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- 3. Apply the "incorrectStyle" (red, centered) look to A36, which now
#        flags a wrong answer. Copy formats only from C10, which already
#        carries that exact style. ---
$ws.Range("C10").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# --- 4. Refresh the summary block (rows 10-12) with the recomputed
#        right/wrong/not-attempted/max counts and marks, keeping every
#        value numeric (fixes the stray "-1" inline-string in C11). ---
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 20
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "19/112"

# --- 5. Fill in the "Correct Ans" values that line up with the newly
#        styled cells. ---
$ws.Range("D16").Value = "Option A"
$ws.Range("A21").Value = "Option C"
$ws.Range("A23").Value = "Option D"
$ws.Range("A28").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A36").Value = "Option B"

# --- 6. Drop the now-unused duplicate "Student Ans / Correct Ans" blocks:
#        the whole G:H columns (rows 15-21) and D:E past row 18
#        (rows 19-40). Clear() removes both value and formatting so the
#        cells disappear entirely and the sheet's used range shrinks back
#        to A5:E40. ---
$ws.Range("G15:H21").Clear() | Out-Null
$ws.Range("D19:E40").Clear() | Out-Null

Write-Host "edit applied"
